# "Update definition of #7"
# Reword the description of project 7 (Rock Paper Scissors) in column D,
# which also updates the H8 concatenation formula's cached value since it
# references D8. Also nudge the layout to match the reflowed text: widen
# column D a bit and shrink row 8's height (the new text is shorter), and
# move the active cell selection to D11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Project 7 description (cell D8, row 8 of the table)
$ws.Range("D8").Value = "On-line version of the popular rock-paper-scissors game. In this version you will be able to just ""watch"" other people play or, if a spot is available, you can play the game against someone else."

# Column D grew a bit wider to accommodate the revised text
$ws.Columns("D").ColumnWidth = 59.67

# Row 8 is shorter now that the text wraps into fewer lines
$ws.Rows(8).RowHeight = 45

# Active selection moved to D11
$ws.Range("D11").Select()
